$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("1er Parcial")
$ws.Range("E5").Value = 36
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 100
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 8.1
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("E6").Value = 21
$ws.Range("F6").Value = 11
$ws.Range("G6").Value = 65.63
$ws.Range("H6").Value = 34.38
$ws.Range("I6").Value = 7.1
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("I7").Value = 7.6
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0

$ws = $wb.Worksheets.Item("2o Parcial")
$ws.Range("E5").Value = 36
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 100
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 7.6
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("E6").Value = 19
$ws.Range("F6").Value = 13
$ws.Range("G6").Value = 59.38
$ws.Range("H6").Value = 40.63
$ws.Range("I6").Value = 8.699999999999999
$ws.Range("J6").Value = 13
$ws.Range("K6").Value = 40.63
$ws.Range("E7").Value = 18
$ws.Range("F7").Value = 7
$ws.Range("G7").Value = 72
$ws.Range("H7").Value = 28
$ws.Range("I7").Value = 8.6
$ws.Range("J7").Value = 7
$ws.Range("K7").Value = 28
$ws.Range("E9").Value = 24
$ws.Range("F9").Value = 9
$ws.Range("G9").Value = 72.73
$ws.Range("H9").Value = 27.27
$ws.Range("J9").Value = 8
$ws.Range("K9").Value = 24.24
$ws.Range("E12").Value = 28
$ws.Range("F12").Value = 9
$ws.Range("G12").Value = 75.68000000000001
$ws.Range("H12").Value = 24.32
$ws.Range("I12").Value = 8.300000000000001
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("E13").Value = 33
$ws.Range("F13").Value = 3
$ws.Range("G13").Value = 91.67
$ws.Range("H13").Value = 8.33
$ws.Range("I13").Value = 8.6
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("E14").Value = 32
$ws.Range("F14").Value = 4
$ws.Range("G14").Value = 88.89
$ws.Range("H14").Value = 11.11
$ws.Range("I14").Value = 8.699999999999999
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("E15").Value = 23
$ws.Range("F15").Value = 12
$ws.Range("G15").Value = 65.70999999999999
$ws.Range("H15").Value = 34.29
$ws.Range("I15").Value = 6.9
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("E16").Value = 23
$ws.Range("F16").Value = 12
$ws.Range("G16").Value = 65.70999999999999
$ws.Range("H16").Value = 34.29
$ws.Range("I16").Value = 6.9
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0

$ws = $wb.Worksheets.Item("3er Parcial")
$ws.Range("E5").Value = 36
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 100
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 8.1
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("E6").Value = 22
$ws.Range("F6").Value = 10
$ws.Range("G6").Value = 68.75
$ws.Range("H6").Value = 31.25
$ws.Range("I6").Value = 7.3
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("I7").Value = 7.6
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("E12").Value = 28
$ws.Range("F12").Value = 9
$ws.Range("G12").Value = 75.68000000000001
$ws.Range("H12").Value = 24.32
$ws.Range("I12").Value = 8.5
$ws.Range("E13").Value = 33
$ws.Range("F13").Value = 3
$ws.Range("G13").Value = 91.67
$ws.Range("H13").Value = 8.33
$ws.Range("I13").Value = 8.699999999999999
$ws.Range("E14").Value = 32
$ws.Range("F14").Value = 4
$ws.Range("G14").Value = 88.89
$ws.Range("H14").Value = 11.11
$ws.Range("I14").Value = 8.699999999999999
$ws.Range("E15").Value = 23
$ws.Range("F15").Value = 12
$ws.Range("G15").Value = 65.70999999999999
$ws.Range("H15").Value = 34.29
$ws.Range("I15").Value = 7.3
$ws.Range("E16").Value = 23
$ws.Range("F16").Value = 12
$ws.Range("G16").Value = 65.70999999999999
$ws.Range("H16").Value = 34.29
$ws.Range("I16").Value = 7.4
